{"js": "// The document contains a paragraph built from three runs:\n//   \"<id>\"  (Courier New, color 7F6000, sz 18)\n//   \"p006r_1\"  (plain run, color 000000)\n//   \"</id>\"  (Courier New, color 7F6000, sz 18)\n// The edit merges these into a single run reading \"<id>p006r_1</id>\"\n// that keeps the formatting of the first (\"<id>\") run.\n\nconst body = context.document.body;\n\n// Locate the opening \"<id>\" tag and the closing \"</id>\" tag.\nconst openResults = body.search(\"<id>\", { matchCase: true });\nconst closeResults = body.search(\"</id>\", { matchCase: true });\nopenResults.load(\"items\");\ncloseResults.load(\"items\");\nawait context.sync();\n\nif (openResults.items.length === 0 || closeResults.items.length === 0) {\n  throw new Error(\"Could not locate <id>...</id> text to merge.\");\n}\n\nconst openRange = openResults.items[0];\nconst closeRange = closeResults.items[0];\n\n// Expand the range to cover \"<id>p006r_1</id>\" in its entirety, then\n// replace it with the same text as a single run. insertText uses the\n// formatting already present at the start of the range being replaced,\n// i.e. the formatting of the original \"<id>\" run, which matches the\n// target output.\nconst fullRange = openRange.expandTo(closeRange);\nfullRange.insertText(\"<id>p006r_1</id>\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The document contains a paragraph built from three runs:\n#   \"<id>\"     (Courier New, color 7F6000, sz 18)\n#   \"p006r_1\"  (plain run, color 000000)\n#   \"</id>\"    (Courier New, color 7F6000, sz 18)\n# The edit merges these into a single run reading \"<id>p006r_1</id>\"\n# that keeps the formatting of the first (\"<id>\") run.\n\n$d = $word.ActiveDocument\n\n# Locate the opening \"<id>\" tag.\n$openRange = $d.Content\n$openFind = $openRange.Find\n$openFind.Text = \"<id>\"\n$openFind.Forward = $true\n$openFind.Wrap = 0\n$openFind.Execute() | Out-Null\n\n# Locate the closing \"</id>\" tag.\n$closeRange = $d.Content\n$closeFind = $closeRange.Find\n$closeFind.Text = \"</id>\"\n$closeFind.Forward = $true\n$closeFind.Wrap = 0\n$closeFind.Execute() | Out-Null\n\n# Remove the old \"p006r_1\" and \"</id>\" runs, leaving only the \"<id>\" run\n# (with its Courier New / 7F6000 formatting) in place.\n$toRemove = $d.Range($openRange.End, $closeRange.End)\n$toRemove.Delete()\n\n# Collapse the selection to the end of the remaining \"<id>\" run and type\n# the rest of the merged text there, so it becomes part of that same run.\n$sel = $word.Selection\n$sel.SetRange($openRange.Start, $openRange.End)\n$sel.Collapse(0)\n$sel.TypeText(\"p006r_1</id>\")\n"}
